# Injuries_Master_Clubs.xlsx update (2025-12-08 run):
#   - snapshot: two newly-flagged SKA injuries are inserted into their
#     alphabetically-sorted slot (Zykov before Korotkiy, Murphy after
#     Korotkiy / before Bikmullin); every other row keeps its data but
#     picks up the fresh scraped_at stamp from this run.
#   - new_injured: the two newly-detected injuries are logged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

# --- bump scraped_at (column K) for the untouched leading rows (2-21) ---
$ws.Cells.Item(2,11).Value  = "2025-12-07T23:09:27.406664+00:00"
$ws.Cells.Item(3,11).Value  = "2025-12-07T23:09:27.406683+00:00"
$ws.Cells.Item(4,11).Value  = "2025-12-07T23:09:27.406693+00:00"
$ws.Cells.Item(5,11).Value  = "2025-12-07T23:09:29.932186+00:00"
$ws.Cells.Item(6,11).Value  = "2025-12-07T23:09:29.932202+00:00"
$ws.Cells.Item(7,11).Value  = "2025-12-07T23:09:32.643213+00:00"
$ws.Cells.Item(8,11).Value  = "2025-12-07T23:09:35.395109+00:00"
$ws.Cells.Item(9,11).Value  = "2025-12-07T23:09:37.771352+00:00"
$ws.Cells.Item(10,11).Value = "2025-12-07T23:09:40.507782+00:00"
$ws.Cells.Item(11,11).Value = "2025-12-07T23:09:45.598718+00:00"
$ws.Cells.Item(12,11).Value = "2025-12-07T23:09:45.598751+00:00"
$ws.Cells.Item(13,11).Value = "2025-12-07T23:09:48.367565+00:00"
$ws.Cells.Item(14,11).Value = "2025-12-07T23:09:55.201191+00:00"
$ws.Cells.Item(15,11).Value = "2025-12-07T23:09:57.561618+00:00"
$ws.Cells.Item(16,11).Value = "2025-12-07T23:10:00.335967+00:00"
$ws.Cells.Item(17,11).Value = "2025-12-07T23:10:00.336000+00:00"
$ws.Cells.Item(18,11).Value = "2025-12-07T23:10:03.083493+00:00"
$ws.Cells.Item(19,11).Value = "2025-12-07T23:10:03.083527+00:00"
$ws.Cells.Item(20,11).Value = "2025-12-07T23:10:03.083549+00:00"
$ws.Cells.Item(21,11).Value = "2025-12-07T23:10:05.391911+00:00"

# --- insert "Зыков Валентин" (СКА) right before "Короткий Матвей" (row 22) ---
$ws.Rows.Item(22).Insert()
$ws.Cells.Item(22,1).Value  = "СКА"
$ws.Cells.Item(22,2).Value  = "СКА"
$ws.Cells.Item(22,3).Value  = "ska"
$ws.Cells.Item(22,4).Value  = "Зыков Валентин"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value  = "90"
$ws.Cells.Item(22,6).Value  = "защитник"
$ws.Cells.Item(22,7).NumberFormat = "@"
$ws.Cells.Item(22,7).Value  = "17992"
$ws.Cells.Item(22,8).Value  = "1369_СКА_зыковвалентин"
$ws.Cells.Item(22,9).Value  = "injured_active"
$ws.Cells.Item(22,10).Value = "https://www.khl.ru/clubs/ska/team/"
$ws.Cells.Item(22,11).Value = "2025-12-07T23:10:05.391939+00:00"

# "Короткий Матвей" is now row 23; refresh its scraped_at.
$ws.Cells.Item(23,11).Value = "2025-12-07T23:10:05.391957+00:00"

# --- insert "Мёрфи Тревор" (СКА) right after "Короткий Матвей" (row 24) ---
$ws.Rows.Item(24).Insert()
$ws.Cells.Item(24,1).Value  = "СКА"
$ws.Cells.Item(24,2).Value  = "СКА"
$ws.Cells.Item(24,3).Value  = "ska"
$ws.Cells.Item(24,4).Value  = "Мёрфи Тревор"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value  = "8"
$ws.Cells.Item(24,6).Value  = "защитник"
$ws.Cells.Item(24,7).NumberFormat = "@"
$ws.Cells.Item(24,7).Value  = "34733"
$ws.Cells.Item(24,8).Value  = "1369_СКА_мерфитревор"
$ws.Cells.Item(24,9).Value  = "injured_active"
$ws.Cells.Item(24,10).Value = "https://www.khl.ru/clubs/ska/team/"
$ws.Cells.Item(24,11).Value = "2025-12-07T23:10:05.391973+00:00"

# --- bump scraped_at for the rows that shifted down by two (old 23-35 -> 25-37) ---
$ws.Cells.Item(25,11).Value = "2025-12-07T23:10:08.225883+00:00"
$ws.Cells.Item(26,11).Value = "2025-12-07T23:10:08.225918+00:00"
$ws.Cells.Item(27,11).Value = "2025-12-07T23:10:08.225940+00:00"
$ws.Cells.Item(28,11).Value = "2025-12-07T23:10:11.017778+00:00"
$ws.Cells.Item(29,11).Value = "2025-12-07T23:10:16.217619+00:00"
$ws.Cells.Item(30,11).Value = "2025-12-07T23:10:16.217652+00:00"
$ws.Cells.Item(31,11).Value = "2025-12-07T23:10:16.217669+00:00"
$ws.Cells.Item(32,11).Value = "2025-12-07T23:10:18.560888+00:00"
$ws.Cells.Item(33,11).Value = "2025-12-07T23:10:18.560915+00:00"
$ws.Cells.Item(34,11).Value = "2025-12-07T23:10:20.909500+00:00"
$ws.Cells.Item(35,11).Value = "2025-12-07T23:10:20.909536+00:00"
$ws.Cells.Item(36,11).Value = "2025-12-07T23:10:23.728001+00:00"
$ws.Cells.Item(37,11).Value = "2025-12-07T23:10:23.728027+00:00"

# --- log the two newly-detected injuries on the new_injured sheet ---
$wsNew = $wb.Worksheets.Item("new_injured")

$wsNew.Cells.Item(2,1).Value = "СКА"
$wsNew.Cells.Item(2,2).Value = "СКА"
$wsNew.Cells.Item(2,3).Value = "Зыков Валентин"
$wsNew.Cells.Item(2,4).Value = "1369_СКА_зыковвалентин"
$wsNew.Cells.Item(2,5).Value = "INJURED_NEW"
$wsNew.Cells.Item(2,6).Value = "2025-12-08T07:10:24.239040+08:00"
$wsNew.Cells.Item(2,7).NumberFormat = "@"
$wsNew.Cells.Item(2,7).Value = "2025-12-08"

$wsNew.Cells.Item(3,1).Value = "СКА"
$wsNew.Cells.Item(3,2).Value = "СКА"
$wsNew.Cells.Item(3,3).Value = "Мёрфи Тревор"
$wsNew.Cells.Item(3,4).Value = "1369_СКА_мерфитревор"
$wsNew.Cells.Item(3,5).Value = "INJURED_NEW"
$wsNew.Cells.Item(3,6).Value = "2025-12-08T07:10:24.239040+08:00"
$wsNew.Cells.Item(3,7).NumberFormat = "@"
$wsNew.Cells.Item(3,7).Value = "2025-12-08"

Write-Output "snapshot + new_injured updated"
